# plantilla_productos.xlsx - "proyecto terminado para presentacion"
#
# For each of the 3 weekday sheets (lunes / martes / jueves) the product
# categories are reorganised:
#   - "plato"  -> "platos"
#   - "bebidas" (3rd category, col G) -> "entrada"
#   - two brand new categories are appended: "postre" (col J) and
#     "bebidas" (col M), each with its own Nombre/Precio/Stock header row
#   - on the "lunes" sheet (the only one with actual stock rows) the old
#     5-row dish list is replaced with a new 4-row list and the now-unused
#     D:I data cells are cleared out.

$wb = $excel.ActiveWorkbook

function Set-CategoryHeaders($ws, $styleRange) {
    # Row 1 category labels (A1 "platos" / D1 "sopa" stay put; G1 becomes
    # "entrada"; J1 "postre" and M1 "bebidas" are brand new columns).
    $ws.Range("A1").Value = "platos"
    $ws.Range("G1").Value = "entrada"
    $ws.Range("J1").Value = "postre"
    $ws.Range("M1").Value = "bebidas"

    # New Nombre/Precio/Stock header triplets for the two new categories.
    $ws.Range("J2").Value = "Nombre"
    $ws.Range("K2").Value = "Precio"
    $ws.Range("L2").Value = "Stock"
    $ws.Range("M2").Value = "Nombre"
    $ws.Range("N2").Value = "Precio"
    $ws.Range("O2").Value = "Stock"

    # Carry over the bold/shaded header formatting used by the sheet
    # (style index differs per sheet) onto the newly added header cells.
    $ws.Range("A1").Copy()
    $ws.Range("J1").PasteSpecial(-4122)
    $ws.Range("A1").Copy()
    $ws.Range("M1").PasteSpecial(-4122)
    $ws.Range("A2:C2").Copy()
    $ws.Range("J2:L2").PasteSpecial(-4122)
    $ws.Range("A2:C2").Copy()
    $ws.Range("M2:O2").PasteSpecial(-4122)
    $excel.CutCopyMode = 0

    # Mirror the existing 3 repeating column-width triplets onto the 2
    # newly used column triplets (J:L, M:O) so the new columns read the
    # same as the first three categories.
    $wA = $ws.Columns.Item(1).ColumnWidth()
    $wB = $ws.Columns.Item(2).ColumnWidth()
    $wC = $ws.Columns.Item(3).ColumnWidth()
    $ws.Columns.Item(10).ColumnWidth = $wA
    $ws.Columns.Item(11).ColumnWidth = $wB
    $ws.Columns.Item(12).ColumnWidth = $wC
    $ws.Columns.Item(13).ColumnWidth = $wA
    $ws.Columns.Item(14).ColumnWidth = $wB
    $ws.Columns.Item(15).ColumnWidth = $wC
}

# ---- lunes (sheet 1) : header relabel + new categories + new dish list ----
$ws1 = $wb.Worksheets.Item("lunes")
Set-CategoryHeaders $ws1 $null

# The "sopa" (D) and old "bebidas" (G, now "entrada") stock rows for the
# dishes that used to live on rows 3-4 are gone - only the "platos" column
# keeps a 4-row dish list now.
$ws1.Range("D3:I4").ClearContents()

$ws1.Range("A3").Value = "tallarines verdes"
$ws1.Range("B3").Value = 6
$ws1.Range("C3").Value = 10

$ws1.Range("A4").Value = "estofado de pollo"
$ws1.Range("B4").Value = 6
$ws1.Range("C4").Value = 10

$ws1.Range("A5").Value = "arroz con pollo"
$ws1.Range("B5").Value = 6
$ws1.Range("C5").Value = 10

$ws1.Range("A6").Value = "chaufa de pollo"
$ws1.Range("B6").Value = 6
$ws1.Range("C6").Value = 10

# Row 7 ("papaliza con aji") no longer exists - the list shrank by one row.
$ws1.Rows.Item(7).Delete()

$ws1.Activate()
$ws1.Range("D6").Select()

# ---- martes / jueves (sheet 2 / 3) : header relabel + new categories only,
#      these sheets never had stock rows beyond the header rows ----
$ws2 = $wb.Worksheets.Item("martes")
Set-CategoryHeaders $ws2 $null

$ws3 = $wb.Worksheets.Item("jueves")
Set-CategoryHeaders $ws3 $null
